$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 231
$ws.Range("I2").Value = 705
$ws.Range("J2").Value = 2823
$ws.Range("K2").Value = 18
$ws.Range("L2").Value = 794
$ws.Range("M2").Value = 47
$ws.Range("N2").Value = 499
$ws.Range("O2").Value = 2
$ws.Range("P2").Value = 13
$ws.Range("Q2").Value = 3
$ws.Range("R2").Value = 35
$ws.Range("S2").Value = 294
$ws.Range("T2").Value = 461
$ws.Range("U2").Value = 21
$ws.Range("V2").Value = 4264
$ws.Range("W2").Value = 1
$ws.Range("X2").Value = 4296
$ws.Range("Y2").Value = 5
$ws.Range("Z2").Value = 56
$ws.Range("AA2").Value = 26
